$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.411.24'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '1.724.68'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.66'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4893'
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2608'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06200'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').Value = '1.721.39'
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07005'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.51'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.571'
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5993'
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.24'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').Value = '26.425.68'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007142'
$ws.Range('E19').Value = '  +2.88%  '
$ws.Range('E20').Value = '  -1.86%  '
$ws.Range('D21').Value = '1.945.37'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.467'
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.586'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.155'
$ws.Range('E24').Value = '  -1.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.18'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.20'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.391'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.87'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.698'
$ws.Range('E29').Value = '  -4.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.947'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07935'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.677'
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04534'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.604'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9945'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6269'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9093'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.394'
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.949'
$ws.Range('E39').Value = '  -6.26%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01481'
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.96'
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.437'
$ws.Range('E43').Value = '  -3.52%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3836'
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.699'
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1155'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05364'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.735'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.08'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.239'
$ws.Range('E50').Value = '  -1.29%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.00'
$ws.Range('E51').Value = '  -0.56%  '
